# Update "想去人数" (interested-count) values in column F
# for worksheets "展览" and "全部类型".
# Both sheets share the same rows/values for this update.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 129
    3  = 410
    4  = 11959
    5  = 1256
    6  = 125
    7  = 25
    9  = 155
    10 = 178
    11 = 263
    13 = 59
    16 = 348
    17 = 1629
    18 = 83
    19 = 920
    20 = 116
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
